# 9th Stab - Cosmetic Changes
# Insert two new "week" columns (Jun_17, Jun_15) in front of the existing
# Jun_13 / Jun_10 columns, shifting the old data two columns to the right,
# and populate the two new columns with the latest analyst-rating rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column B. The previous B ("Jun_13")
# and C ("Jun_10") columns - and all their data/styles - shift to D and E.
$ws.Columns("B:C").Insert()

# --- Header row -----------------------------------------------------------
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# --- Fill the two new columns with the default "UN" (unchanged) marker ----
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# --- New rating-change callouts -------------------------------------------
# Zacks Investment Research (row 5): 6/16/2018 Sell -> Hold upgrade, new column B
$ws.Range("B5").Value = "6/16/2018,Upgrades,Sell -> Hold,"
$ws.Range("B5").Interior.Color = 13434828

# ValuEngine (row 11): 6/15/2018 Strong Sell -> Sell upgrade, new column C
$ws.Range("C11").Value = "6/15/2018,Upgrades,Strong Sell -> Sell,"
$ws.Range("C11").Interior.Color = 13434828

# --- Column widths ----------------------------------------------------------
# Column B (new "Jun_17") keeps the workbook default width - no change needed.
# Columns C ("Jun_15", new), D ("Jun_13", was B) and E ("Jun_10", was C) all
# carry the same explicit 8.0-unit width the original "Jun_10" column had.
$ws.Columns.Item(3).ColumnWidth = 7.1666666667
$ws.Columns.Item(4).ColumnWidth = 7.1666666667
$ws.Columns.Item(5).ColumnWidth = 7.1666666667
